# Update the "raices multiples" results table:
#  - Row 2 (iteration 1) gets corrected values for the specific case.
#  - Row 3 (iteration 2) gets corrected values and becomes the last data row.
#  - Rows 4-8 (iterations 3-7) are removed entirely, shrinking the table
#    from A1:D8 down to A1:D3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Leading apostrophe forces these numeric-looking values to be stored as
# text, matching the original table's text-based number formatting
# (e.g. "2.0" / "0.0" rather than numeric 2 / 0).
$ws.Range("B2").Value = "'1.5"
$ws.Range("C2").Value = "'0.25"
$ws.Range("D2").Value = "'0.5"

$ws.Range("B3").Value = "'2.0"
$ws.Range("C3").Value = "'0.0"
$ws.Range("D3").Value = "'0.0"

# The apostrophe prefix also marks the cells with a "quote prefix" style;
# reset back to the plain Normal style so no stray formatting is left
# behind (the source cells were unstyled).
$ws.Range("B2:D3").Style = "Normal"

# Remove the now-obsolete iteration rows 3-7 (sheet rows 4-8).
$ws.Range("A4:D8").Clear()
